$wb = $excel.ActiveWorkbook
$ws1 = $wb.Sheets.Item("Static")
$ws2 = $wb.Sheets.Item("Dynamic")

# ---------------------------------------------------------------------------
# Dynamic sheet (sheet2): fill in the Name_Change / Size_Change values that
# the bot captured, and add a second column width + new selection.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Name_Change"
$ws2.Range("B2").Value = "20:00, 18 March 2023"
$ws2.Range("A3").Value = "Size_Change"
$ws2.Range("B3").Value = "144,909 bytes"
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Static sheet (sheet1): existing header/URL/size rows get re-centered, then
# a whole new "mail + credentials" configuration block is appended below.
# ---------------------------------------------------------------------------

# Row 1 - header row (unchanged text, new center+middle alignment style)
$ws1.Range("A1").Value = "Name"
$ws1.Range("B1").Value = "Value"
$ws1.Range("C1").Value = "Description"

# Row 2 - URL row
$ws1.Range("A2").Value = "URL"
$ws1.Range("B2").Value = "https://en.wikipedia.org/w/index.php?title=A._P._J._Abdul_Kalam&action=history"
$ws1.Range("C2").Value = "URL to work"

# Row 3 - Prefered_Size row
$ws1.Range("A3").Value = "Prefered_Size"
$ws1.Range("B3").Value = 500
$ws1.Range("C3").Value = "Provide the value in bytes"

# Re-style rows 1:3 -> horizontal=center, vertical=center (new style #2)
$ws1.Range("A1:C3").Style = "Normal"
$ws1.Range("A1:C3").HorizontalAlignment = -4108
$ws1.Range("A1:C3").VerticalAlignment = -4108

# Row 4 - Mail Body (long wrapped text, centered + wrap, taller row)
$ws1.Range("A4").Style = "Normal"
$ws1.Range("A4").Value = "Mail Body"
$ws1.Range("A4").HorizontalAlignment = -4108
$ws1.Range("A4").VerticalAlignment = -4108

$ws1.Range("B4").Style = "Normal"
$ws1.Range("B4").Value = "Hi Sir," + [char]10 + "There are large changes on Website" + [char]10 + "Sent From Bot"
$ws1.Range("B4").HorizontalAlignment = -4108
$ws1.Range("B4").VerticalAlignment = -4108
$ws1.Range("B4").WrapText = $true

$ws1.Range("C4").Style = "Normal"
$ws1.Range("C4").Value = "What you want to write in mail body"
$ws1.Range("C4").HorizontalAlignment = -4108
$ws1.Range("C4").VerticalAlignment = -4108

$ws1.Rows.Item(4).RowHeight = 60

# Row 5 - Mail Subject
$ws1.Range("A5").Style = "Normal"
$ws1.Range("A5").Value = "Mail Subject"
$ws1.Range("A5").HorizontalAlignment = -4108
$ws1.Range("A5").VerticalAlignment = -4108

$ws1.Range("B5").Style = "Normal"
$ws1.Range("B5").Value = "Major Changes on Site"
$ws1.Range("B5").HorizontalAlignment = -4108
$ws1.Range("B5").VerticalAlignment = -4108

$ws1.Range("C5").Style = "Normal"
$ws1.Range("C5").Value = "Subject of Mail"
$ws1.Range("C5").HorizontalAlignment = -4108
$ws1.Range("C5").VerticalAlignment = -4108

# Row 6 - Send To
$ws1.Range("A6").Style = "Normal"
$ws1.Range("A6").Value = "Send To"
$ws1.Range("A6").HorizontalAlignment = -4108
$ws1.Range("A6").VerticalAlignment = -4108

$ws1.Range("B6").Style = "Normal"
$ws1.Range("B6").Value = "gauravkeny1@gmail.com"
$ws1.Range("B6").HorizontalAlignment = -4108
$ws1.Range("B6").Interior.ColorIndex = 6
$ws1.Range("B6").Interior.ColorIndex = -4142

$ws1.Range("C6").Style = "Normal"
$ws1.Range("C6").Value = "Enter multiple email id in comma seperated form"
$ws1.Range("C6").HorizontalAlignment = -4108
$ws1.Range("C6").VerticalAlignment = -4108

# Row 7 - Send CC (value column left blank on purpose)
$ws1.Range("A7").Style = "Normal"
$ws1.Range("A7").Value = "Send CC"
$ws1.Range("A7").HorizontalAlignment = -4108
$ws1.Range("A7").VerticalAlignment = -4108

$ws1.Range("B7").Style = "Normal"
$ws1.Range("B7").HorizontalAlignment = -4108
$ws1.Range("B7").VerticalAlignment = -4108

$ws1.Range("C7").Style = "Normal"
$ws1.Range("C7").Value = "Email Ids for cc Part"
$ws1.Range("C7").HorizontalAlignment = -4108
$ws1.Range("C7").VerticalAlignment = -4108

# Row 8 - Send From
$ws1.Range("A8").Style = "Normal"
$ws1.Range("A8").Value = "Send From"
$ws1.Range("A8").HorizontalAlignment = -4108
$ws1.Range("A8").VerticalAlignment = -4108

$ws1.Range("B8").Style = "Normal"
$ws1.Range("B8").Value = "gauravkeny1@gmail.com"
$ws1.Range("B8").HorizontalAlignment = -4108
$ws1.Range("B8").Interior.ColorIndex = 6
$ws1.Range("B8").Interior.ColorIndex = -4142

$ws1.Range("C8").Style = "Normal"
$ws1.Range("C8").Value = "From which email Id need to send mail"
$ws1.Range("C8").HorizontalAlignment = -4108
$ws1.Range("C8").VerticalAlignment = -4108

# Row 10 - Enc Client ID (row 9 intentionally left blank)
$ws1.Range("A10").Style = "Normal"
$ws1.Range("A10").Value = "Enc Client ID"
$ws1.Range("A10").HorizontalAlignment = -4108
$ws1.Range("A10").VerticalAlignment = -4108

$ws1.Range("B10").Style = "Normal"
$ws1.Range("B10").Value = "lRMNNqzor31QjVFoWEm+oI5JpLHDEhVBNROJe2hWdd6saS2fitYM2zAL19uO5Hv9Uum30eANlraiO/P+Gkv9hxcawzQ8oWdHZFb9xwK0y6rKgexShlbQEWkk3/ZgBahb+JC2dMuHmnVJKqqeFA=="
$ws1.Range("B10").WrapText = $true

$ws1.Range("C10").Style = "Normal"
$ws1.Range("C10").Value = "Encrypted Client ID"
$ws1.Range("C10").HorizontalAlignment = -4108
$ws1.Range("C10").VerticalAlignment = -4108

$ws1.Rows.Item(10).RowHeight = 45

# Row 11 - Enc Client Secret
$ws1.Range("A11").Style = "Normal"
$ws1.Range("A11").Value = "Enc Client Secret"
$ws1.Range("A11").HorizontalAlignment = -4108
$ws1.Range("A11").VerticalAlignment = -4108

$ws1.Range("B11").Style = "Normal"
$ws1.Range("B11").Value = "S+Rq0hUcXnUAFSA47AfGpKVU9uhkW1PFMwlM/cpzHvN6yCYdplU8L1f6xnBGyWmsQHXlZ0mu7+Us6M0zuNJiEi42xT1kuDE="
$ws1.Range("B11").WrapText = $true

$ws1.Range("C11").Style = "Normal"
$ws1.Range("C11").Value = "Encrypted Client ID"
$ws1.Range("C11").HorizontalAlignment = -4108
$ws1.Range("C11").VerticalAlignment = -4108

$ws1.Rows.Item(11).RowHeight = 30

# Column widths (auto-size columns A and C for the new, wider content)
$ws1.Columns.Item(1).AutoFit() | Out-Null
$ws1.Columns.Item(3).AutoFit() | Out-Null

# Final selection / active sheet state
$ws1.Activate() | Out-Null
$ws1.Range("C11").Select() | Out-Null
